$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62 - this shifts the existing row 62..126
# (and everything below) down by one row, extending the data through row 127.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly price-report record.
$ws.Cells.Item(62, 1).Value = 7
$ws.Cells.Item(62, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(62, 3).Value = "Ñuble"
$ws.Cells.Item(62, 4).Value = 44601
$ws.Cells.Item(62, 5).Value = 16
$ws.Cells.Item(62, 6).Value = 100112045
$ws.Cells.Item(62, 7).Value = "Zapallo"
$ws.Cells.Item(62, 8).Value = "Camote"
$ws.Cells.Item(62, 9).Value = "1a (cosecha)"
$ws.Cells.Item(62, 10).Value = 240
$ws.Cells.Item(62, 11).Value = 350
$ws.Cells.Item(62, 12).Value = 400
$ws.Cells.Item(62, 13).Value = 375
$ws.Cells.Item(62, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(62, 15).Value = "Región del Maule"
$ws.Cells.Item(62, 16).Value = 375
$ws.Cells.Item(62, 17).Value = 1
$ws.Cells.Item(62, 18).Value = "Hortaliza"
